# Config.xlsx: "Change credential management and storage."
# Replace the plaintext Email/Password settings rows with a single row that
# references the Windows Credential Manager asset used to store the
# ACME system credential.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Email" row (row 15) in place to hold the new credential-asset
# name / value pair.
$ws.Range("A15").Value = "ACME_CredentialAssetNameInWindowsCredentialsManager"
$ws.Range("B15").Value = "ACMESystem1_Credential"

# Give B15 the same "data value" styling used by the rest of column B in this
# block (e.g. B14), since it previously held an unstyled literal string.
$ws.Range("B15").Font.Bold = $false

# Remove the old "Password" row (row 16) entirely - its content has been
# superseded by the single credential-asset row above, and every row below
# shifts up by one.
$ws.Rows("16:16").Delete()

# Resize column A to fit the new, longer setting name.
$ws.Columns("A:A").ColumnWidth = 54.675

# Update the view: scroll so row 5 is at the top and select A9:C10 (the
# "Initialization Constants" header block).
$ws.Activate()
$ws.Range("A9:C10").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
